$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: corrected hydrogen demand value
$ws.Range("B3").Value = 7573949.92189134

# D3: value removed -> cell becomes blank
$ws.Range("D3").ClearContents()

# C4: corrected methanol chemicals value
$ws.Range("C4").Value = 220.0670067352325

# C5: corrected ammonia chemicals value
$ws.Range("C5").Value = 6343.984116856707

# Row 7 ("Other") renamed to "Biogas", D7 value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 12838.86571380623

# New row 8: "Other" row, carrying forward the label style from column A
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 9593.717256879121
